$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 8335304
$ws.Range("I19").Value = 15625781
$ws.Range("J19").Value = 3329.4285
$ws.Range("K19").Value = 15625781
$ws.Range("L19").Value = 3329.4285
$ws.Range("M19").Value = -15625606
$ws.Range("N19").Value = -3679.4285

# Row 39
$ws.Range("H39").Value = 35715050
$ws.Range("I39").Value = 44.666668
$ws.Range("J39").Value = 100002050
$ws.Range("K39").Value = 134.000004
$ws.Range("L39").Value = 300006150
$ws.Range("M39").Value = 161.999996
$ws.Range("N39").Value = -300006742

# Row 116
$ws.Range("H116").Value = 16745563
$ws.Range("I116").Value = 8335901
$ws.Range("J116").Value = 27958446
$ws.Range("K116").Value = 8335901
$ws.Range("L116").Value = 27958446
$ws.Range("M116").Value = -8332459
$ws.Range("N116").Value = -27965330

# Row 137
$ws.Range("H137").Value = 19244528
$ws.Range("I137").Value = 4465533.5
$ws.Range("J137").Value = 56863790
$ws.Range("K137").Value = 13396600.5
$ws.Range("L137").Value = 170591370
$ws.Range("M137").Value = -13394050.5
$ws.Range("N137").Value = -170596470

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 5026058
$ws.Range("I61").Value = 2605711.2
$ws.Range("J61").Value = 14707444
$ws.Range("K61").Value = 2605711.2
$ws.Range("L61").Value = 14707444
$ws.Range("M61").Value = -2605499.2
$ws.Range("N61").Value = -14707868

# Row 74
$ws.Range("H74").Value = 40352400
$ws.Range("I74").Value = 38462360
$ws.Range("J74").Value = 44447480
$ws.Range("K74").Value = 38462360
$ws.Range("L74").Value = 44447480
$ws.Range("M74").Value = -38461486
$ws.Range("N74").Value = -44449228

# Row 77
$ws.Range("H77").Value = 40352400
$ws.Range("I77").Value = 38462360
$ws.Range("J77").Value = 44447480
$ws.Range("K77").Value = 192311800
$ws.Range("L77").Value = 222237400
$ws.Range("M77").Value = -192307432
$ws.Range("N77").Value = -222246136

# Row 132
$ws.Range("H132").Value = 17682850
$ws.Range("I132").Value = 22230012
$ws.Range("J132").Value = 7938934.5
$ws.Range("K132").Value = 66690036
$ws.Range("L132").Value = 23816803.5
$ws.Range("M132").Value = -66687506
$ws.Range("N132").Value = -23821863.5

# Row 136
$ws.Range("H136").Value = 5026058
$ws.Range("I136").Value = 2605711.2
$ws.Range("J136").Value = 14707444
$ws.Range("K136").Value = 7817133.600000001
$ws.Range("L136").Value = 44122332
$ws.Range("M136").Value = -7814583.600000001
$ws.Range("N136").Value = -44127432

$ws = $wb.Worksheets.Item("BSM")
# Row 57
$ws.Range("H57").Value = 20780
$ws.Range("J57").Value = 20780
$ws.Range("L57").Value = 20780
$ws.Range("N57").Value = -22220

# Row 132
$ws.Range("H132").Value = 29610.666
$ws.Range("J132").Value = 29610.666
$ws.Range("L132").Value = 29610.666
$ws.Range("N132").Value = -39730.666

# Row 133
$ws.Range("H133").Value = 40780
$ws.Range("J133").Value = 40780
$ws.Range("L133").Value = 40780
$ws.Range("N133").Value = -50900

# Row 134
$ws.Range("H134").Value = 35715640
$ws.Range("I134").Value = 50000932
$ws.Range("J134").Value = 7145056.5
$ws.Range("K134").Value = 150002796
$ws.Range("L134").Value = 21435169.5
$ws.Range("M134").Value = -150000261
$ws.Range("N134").Value = -21440239.5

# Row 136
$ws.Range("H136").Value = 20780
$ws.Range("J136").Value = 20780
$ws.Range("L136").Value = 20780
$ws.Range("N136").Value = -30980

# Row 137
$ws.Range("H137").Value = 55790
$ws.Range("J137").Value = 55790
$ws.Range("L137").Value = 55790
$ws.Range("N137").Value = -65990

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# Row 139
$ws.Range("H139").Value = 72340.336
$ws.Range("I139").Value = 35000
$ws.Range("J139").Value = 147021
$ws.Range("K139").Value = 35000
$ws.Range("L139").Value = 147021
$ws.Range("N139").Value = -157301
$ws.Range("M139").Value = -29860

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# Row 141
$ws.Range("H141").Value = 100840
$ws.Range("J141").Value = 100840
$ws.Range("L141").Value = 100840
$ws.Range("N141").Value = -111200

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1604747
$ws.Range("I31").Value = 1592
$ws.Range("J31").Value = 5211846
$ws.Range("K31").Value = 1592
$ws.Range("L31").Value = 5211846
$ws.Range("M31").Value = -1297
$ws.Range("N31").Value = -5212436

# Row 34
$ws.Range("H34").Value = 1604747
$ws.Range("I34").Value = 1592
$ws.Range("J34").Value = 5211846
$ws.Range("K34").Value = 1592
$ws.Range("L34").Value = 5211846
$ws.Range("M34").Value = -1390
$ws.Range("N34").Value = -5212250

# Row 107
$ws.Range("H107").Value = 719.3214
$ws.Range("I107").Value = 296
$ws.Range("J107").Value = 860.4286
$ws.Range("K107").Value = 296
$ws.Range("L107").Value = 860.4286
$ws.Range("M107").Value = 1624
$ws.Range("N107").Value = -4700.4286

# Row 122
$ws.Range("H122").Value = 13080
$ws.Range("I122").Value = 13080
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 39240
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -36790
$ws.Range("N122").ClearContents()

# Row 134
$ws.Range("H134").Value = 1255122.1
$ws.Range("I134").Value = 4796.1113
$ws.Range("J134").Value = 8006883
$ws.Range("K134").Value = 14388.3339
$ws.Range("L134").Value = 24020649
$ws.Range("M134").Value = -11853.3339
$ws.Range("N134").Value = -24025719

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 281.2143
$ws.Range("I8").Value = 281.2143
$ws.Range("K8").Value = 843.6428999999999
$ws.Range("M8").Value = -704.6428999999999

# Row 36
$ws.Range("H36").Value = 142894990
$ws.Range("I36").Value = 366.66666
$ws.Range("J36").Value = 250065970
$ws.Range("K36").Value = 1099.99998
$ws.Range("L36").Value = 750197910
$ws.Range("M36").Value = -930.9999800000001
$ws.Range("N36").Value = -750198248

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 12539457
$ws.Range("I132").Value = 12381848
$ws.Range("J132").Value = 12989769
$ws.Range("K132").Value = 37145544
$ws.Range("L132").Value = 38969307
$ws.Range("M132").Value = -37143014
$ws.Range("N132").Value = -38974367

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 8311
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 10024.556
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 10024.556
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -10614.556

# Row 27
$ws.Range("H27").Value = 8311
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 10024.556
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 10024.556
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -10238.556

# Row 132
$ws.Range("H132").Value = 17861004
$ws.Range("I132").Value = 35715760
$ws.Range("J132").Value = 6248.25
$ws.Range("K132").Value = 107147280
$ws.Range("L132").Value = 18744.75
$ws.Range("M132").Value = -107144750
$ws.Range("N132").Value = -23804.75

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 7013.7295
$ws.Range("I136").Value = 5844.304
$ws.Range("J136").Value = 8934.929
$ws.Range("K136").Value = 17532.912
$ws.Range("L136").Value = 26804.787
$ws.Range("M136").Value = -14982.912
$ws.Range("N136").Value = -31904.787
